$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Normalise the date formatting on the two rows that were appended by a
#    different tool (rows 6 & 7) so they share the same "yyyy-mm-dd" number
#    format as the rest of the Date Of Reservation column.
# ---------------------------------------------------------------------------
$ws.Range("C6").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C7").NumberFormat = "yyyy\-mm\-dd"

# ---------------------------------------------------------------------------
# 2) New headers for the two extra columns.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "Table Type"
$ws.Range("F1").Value = "Table Count"

# ---------------------------------------------------------------------------
# 3) Fill in "Table Type" / "Table Count" for the existing reservation rows.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "2 seat"
$ws.Range("F2").Value = 2

$ws.Range("E6").Value = "4 seat"
$ws.Range("F6").Value = 1

$ws.Range("E3").Value = "10 seat"
$ws.Range("F3").Value = 1

$ws.Range("E4").Value = "2 seat"
$ws.Range("F4").Value = 1

$ws.Range("E5").Value = "2 seat"
$ws.Range("F5").Value = 3

$ws.Range("E7").Value = "4 seat"
$ws.Range("F7").Value = 1

# ---------------------------------------------------------------------------
# 4) New reservation row 8, including its own Table Type / Table Count.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 45410
$ws.Range("C8").NumberFormat = "yyyy-MM-dd"
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = "8 seat"
$ws.Range("F8").Value = 1

# ---------------------------------------------------------------------------
# 5) Cosmetics: column widths for the new columns, page orientation and the
#    active selection, matching where the authors left the cursor.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 9.9
$ws.Columns.Item(6).ColumnWidth = 10.96

$ws.PageSetup.Orientation = 1

$ws.Range("F8").Select()
